$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.053.29"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "  -2.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.819.23"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "  -1.41%  "

$ws.Range("E4").Value = "  -1.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.61"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  -2.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  -1.08%  "

$ws.Range("E7").Value = "  -2.10%  "

$ws.Range("E8").Value = "  -1.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07205"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "  -1.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8420"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  -3.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.84"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "  -3.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.822.20"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "  -1.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.618"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "  -1.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07072"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  -0.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.269"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "  -3.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.99"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "  +1.14%  "

$ws.Range("E17").Value = "  -1.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008796"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "  -1.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  -1.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.95"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "  -3.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.104.88"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  -2.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.111"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "  -1.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.80"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "  -2.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.042.61"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = "  -1.87%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.975"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = "  -1.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.71"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "  -2.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.233"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = "  +4.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.30"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = "  -1.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.200"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = "  -3.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.89"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = "  -2.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08790"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = "  -1.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.175"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  -4.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.987"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = "  +2.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7381"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  -4.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.407"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = "  -3.03%  "

$ws.Range("E36").Value = "  -1.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.099"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  -3.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01958"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = "  -0.68%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05232"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = "  -1.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.254"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.876"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "  -0.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1687"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5019"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "  -1.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.555"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "  -2.25%  "

$ws.Range("E45").Value = "  -1.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4734"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("E47").Value = "  -3.19%  "

$ws.Range("E48").Value = "  -1.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06363"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "  -1.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.645"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = "  -2.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.875"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = "  +1.87%  "
